$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Worksheet")

# Grab the picture that sits over the header area before the row shift so we
# can nudge it back up with the deleted row (Excel moves "move with cells"
# anchored shapes automatically; this engine needs an explicit push).
$pic = $ws.Shapes.Item(1)
$picTop = $pic.Top

# The new earnings release (Q4 24, reported 02/13/2025) is removed - delete
# its row outright so everything below ripples up by one row.
$ws.Rows(2).Delete()

# Keep the logo anchored to the same relative header position now that row 2
# is gone (shift up by one default row height).
$pic.Top = $picTop - 15

# Reflect where the author's cursor ended up after the edit.
$ws.Activate()
$ws.Range("C8").Select()
